$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "金风科技"
$ws.Range("A3").Value = "航天发展"
$ws.Range("A4").Value = "中国卫星"
$ws.Range("B4").Value = "蓝色光标"
$ws.Range("C4").Value = "利欧股份"
$ws.Range("A5").Value = "雷科防务"
$ws.Range("B5").Value = "航天电子"
$ws.Range("C5").Value = "平潭发展"
$ws.Range("A6").Value = "三花智控"
$ws.Range("C6").Value = "雷科防务"
$ws.Range("A7").Value = "顺灏股份"
$ws.Range("B7").Value = "雷科防务"
$ws.Range("C7").Value = "三花智控"
$ws.Range("A8").Value = "利欧股份"
$ws.Range("B8").Value = "三花智控"
$ws.Range("C8").Value = "神剑股份"
$ws.Range("A9").Value = "蓝色光标"
$ws.Range("B9").Value = "航天发展"
$ws.Range("C9").Value = "泰尔股份"
$ws.Range("A10").Value = "航天电子"
$ws.Range("B10").Value = "中国卫通"
$ws.Range("C10").Value = "中国卫星"
$ws.Range("A11").Value = "神剑股份"
$ws.Range("B11").Value = "岩山科技"
$ws.Range("C11").Value = "蓝色光标"
$ws.Range("A12").Value = "中国卫通"
$ws.Range("B12").Value = "五洲新春"
$ws.Range("C12").Value = "山子高科"
$ws.Range("A13").Value = "泰尔股份"
$ws.Range("B13").Value = "顺灏股份"
$ws.Range("C13").Value = "航天电子"
$ws.Range("A14").Value = "王子新材"
$ws.Range("A15").Value = "岩山科技"
$ws.Range("B15").Value = "信维通信"
$ws.Range("C15").Value = "顺灏股份"
$ws.Range("A16").Value = "五洲新春"
$ws.Range("B16").Value = "泰尔股份"
$ws.Range("C16").Value = "大业股份"
$ws.Range("A17").Value = "联创光电"
$ws.Range("B17").Value = "王子新材"
$ws.Range("C17").Value = "中国卫通"
$ws.Range("A18").Value = "万向钱潮"
$ws.Range("B18").Value = "北斗星通"
$ws.Range("A19").Value = "山子高科"
$ws.Range("B19").Value = "创新医疗"
$ws.Range("C19").Value = "五洲新春"
$ws.Range("A20").Value = "平潭发展"
$ws.Range("B20").Value = "烽火通信"
$ws.Range("C20").Value = "岩山科技"
$ws.Range("A21").Value = "北斗星通"
$ws.Range("B21").Value = "万向钱潮"
